# Aggiorno file need_to_buy.xlsx da R
# Refresh the forecast/buy table: the rolling date window advances by one
# day (A2:A15) and the corresponding TISG/fcs/buy/MYDIR/need_to_buy values
# (B2:F15) are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46037, 12405.4205184334, 11731.4651981752, 20475.86,  8349.10819295314, -16.4702753696512),
    @(46038, 12168.9196322615, 10842.5123607866, 12075.86,  8273.67576886036, 293.347005401957),
    @(46039, 4640.79444033817, 7347.83089292821, 12075.86,  7753.8260990872,  126.074874667309),
    @(46040, 4533.59182166025, 7265.25119703893, 12075.86,  7719.15386334973, 121.189377516194),
    @(46041, 12011.2256033508, 11271.8156368244, 12075.86,  8099.12789137558, 303.961813674999),
    @(46042, 12733.0845017406, 12107.3005556413, 12075.86,  8656.54013840321, 361.999195585189),
    @(46043, 12733.0845017406, 12373.2046842254, 12075.86,  8656.54013840321, 373.07853427619),
    @(46044, 12733.0845017406, 12376.0404158932, 12075.86,  8656.54013840321, 373.19668976235),
    @(46045, 12733.0845017406, 11677.6456954357, 12075.86,  8656.54013840321, 344.096909743286),
    @(46046, 5107.54318705847, 8311.07642608366, 12075.86,  8263.39417905842, 187.442108547587),
    @(46047, 5000.01932310789, 8366.61664849186, 12075.86,  8255.64962845675, 189.433594872859),
    @(46048, 12315.9682835607, 12482.0582361346, 12075.86,  8313.23194980497, 363.309591080817),
    @(46049, 12315.9682835607, 13078.5188697831, 12075.86,  8313.23194980497, 388.162117482837),
    @(46050, 12315.9682835607, 13358.4460986358, 12075.86,  8313.23194980497, 399.825752018366)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
